# Update premium values and append a new tracking row on CRSRCoveredCalls,
# then leave that sheet as the active tab (as it was when the workbook was
# last saved by the author).
$wb = $excel.ActiveWorkbook

$wsCalls = $wb.Worksheets.Item("CRSRCoveredCalls")

# Existing rows whose premium-paid values changed
$wsCalls.Range("B6").Value = 15
$wsCalls.Range("B7").Value = 15

# New tracking row appended at the bottom of the table
$wsCalls.Range("A14").Value = 62
$wsCalls.Range("B14").Value = 0

# Selection stays on B14, and CRSRCoveredCalls becomes the active sheet
$wsCalls.Range("B14").Select()
$wsCalls.Activate()
